$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Coin name / Link / Volume text cells (plain text, safe to assign directly)
$ws.Range("E2").Value = "  +7.27%  "
$ws.Range("E3").Value = "  +4.64%  "
$ws.Range("E4").Value = "  +0.60%  "
$ws.Range("E5").Value = "  +10.25%  "
$ws.Range("E6").Value = "  +3.27%  "
$ws.Range("E7").Value = "  +3.55%  "
$ws.Range("E8").Value = "  +0.12%  "
$ws.Range("E9").Value = "  +4.67%  "
$ws.Range("E10").Value = "  +11.64%  "
$ws.Range("E11").Value = "  +4.26%  "
$ws.Range("E12").Value = "  +6.53%  "
$ws.Range("E13").Value = "  +1.99%  "
$ws.Range("E14").Value = "  +5.37%  "
$ws.Range("E15").Value = "  +5.43%  "
$ws.Range("E16").Value = "  +4.65%  "
$ws.Range("E17").Value = "  +4.83%  "
$ws.Range("E18").Value = "  +7.16%  "
$ws.Range("E19").Value = "  +5.69%  "
$ws.Range("E20").Value = "  +4.70%  "
$ws.Range("E21").Value = "  +1.07%  "
$ws.Range("E23").Value = "  +5.49%  "
$ws.Range("E24").Value = "  +4.33%  "
$ws.Range("E25").Value = "  +10.17%  "
$ws.Range("E26").Value = "  -0.74%  "
$ws.Range("E27").Value = "  +6.65%  "
$ws.Range("E28").Value = "  +7.38%  "
$ws.Range("E29").Value = "  +2.07%  "
$ws.Range("E30").Value = "  +10.62%  "
$ws.Range("E31").Value = "  +4.10%  "
$ws.Range("E32").Value = "  +12.61%  "
$ws.Range("E33").Value = "  +3.84%  "
$ws.Range("E34").Value = "  +15.84%  "
$ws.Range("E35").Value = "  +2.89%  "
$ws.Range("B36").Value = "RenderToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("E36").Value = "  +11.38%  "
$ws.Range("B37").Value = "Kaspa"
$ws.Range("C37").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("E37").Value = "  +8.78%  "
$ws.Range("E38").Value = "  +13.07%  "
$ws.Range("E39").Value = "  +6.03%  "
$ws.Range("E40").Value = "  +11.45%  "
$ws.Range("E41").Value = "  +11.09%  "
$ws.Range("E42").Value = "  +6.38%  "
$ws.Range("E43").Value = "  +8.19%  "
$ws.Range("B44").Value = "MultiversX"
$ws.Range("C44").Value = "https://coinranking.com/coin/omwkOTglq+multiversx-egld"
$ws.Range("E44").Value = "  +3.84%  "
$ws.Range("B45").Value = "Celestia"
$ws.Range("C45").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("E45").Value = "  +11.77%  "
$ws.Range("E46").Value = "  +0.19%  "
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("E47").Value = "  +13.63%  "
$ws.Range("B48").Value = "Aave"
$ws.Range("C48").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("E48").Value = "  +6.52%  "
$ws.Range("E49").Value = "  +21.12%  "
$ws.Range("E51").Value = "  +18.63%  "

# Update Price cells: use a leading apostrophe to force text interpretation
# (prevents Excel auto-converting numeric-looking strings to numbers/dates),
# then reset the cell style to Normal so no stray quote-prefix style lingers.
$c = $ws.Range("D2")
$c.Value = "'45.410.09"
$c.Style = "Normal"
$c = $ws.Range("D3")
$c.Value = "'2.382.20"
$c.Style = "Normal"
$c = $ws.Range("D5")
$c.Value = "'113.15"
$c.Style = "Normal"
$c = $ws.Range("D6")
$c.Value = "'317.53"
$c.Style = "Normal"
$c = $ws.Range("D9")
$c.Value = "'0.626"
$c.Style = "Normal"
$c = $ws.Range("D10")
$c.Value = "'42.69"
$c.Style = "Normal"
$c = $ws.Range("D11")
$c.Value = "'0.0931"
$c.Style = "Normal"
$c = $ws.Range("D15")
$c.Value = "'15.87"
$c.Style = "Normal"
$c = $ws.Range("D16")
$c.Value = "'2.743.84"
$c.Style = "Normal"
$c = $ws.Range("D17")
$c.Value = "'2.384.95"
$c.Style = "Normal"
$c = $ws.Range("D18")
$c.Value = "'45.317.95"
$c.Style = "Normal"
$c = $ws.Range("D19")
$c.Value = "'7.63"
$c.Style = "Normal"
$c = $ws.Range("D21")
$c.Value = "'13.45"
$c.Style = "Normal"
$c = $ws.Range("D22")
$c.Value = "'74.76"
$c.Style = "Normal"
$c = $ws.Range("D23")
$c.Value = "'3.55"
$c.Style = "Normal"
$c = $ws.Range("D24")
$c.Value = "'269.33"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.Value = "'11.31"
$c.Style = "Normal"
$c = $ws.Range("D28")
$c.Value = "'7.54"
$c.Style = "Normal"
$c = $ws.Range("D30")
$c.Value = "'39.41"
$c.Style = "Normal"
$c = $ws.Range("D31")
$c.Value = "'22.93"
$c.Style = "Normal"
$c = $ws.Range("D32")
$c.Value = "'0.0952"
$c.Style = "Normal"
$c = $ws.Range("D33")
$c.Value = "'170.26"
$c.Style = "Normal"
$c = $ws.Range("D36")
$c.Value = "'4.97"
$c.Style = "Normal"
$c = $ws.Range("D37")
$c.Value = "'0.120"
$c.Style = "Normal"
$c = $ws.Range("D38")
$c.Value = "'3.05"
$c.Style = "Normal"
$c = $ws.Range("D39")
$c.Value = "'0.0365"
$c.Style = "Normal"
$c = $ws.Range("D40")
$c.Value = "'4.00"
$c.Style = "Normal"
$c = $ws.Range("D42")
$c.Value = "'105.43"
$c.Style = "Normal"
$c = $ws.Range("D44")
$c.Value = "'71.49"
$c.Style = "Normal"
$c = $ws.Range("D45")
$c.Value = "'13.35"
$c.Style = "Normal"
$c = $ws.Range("D47")
$c.Value = "'5.79"
$c.Style = "Normal"
$c = $ws.Range("D48")
$c.Value = "'116.72"
$c.Style = "Normal"
$c = $ws.Range("D49")
$c.Value = "'1.67"
$c.Style = "Normal"
$c = $ws.Range("D50")
$c.Value = "'9.35"
$c.Style = "Normal"
